$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-20: update parameter list
$ws.Cells.Item(2, 1).Value = "Number of 463L master pallet crates"
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = "-"

$ws.Cells.Item(3, 1).Value = "Number of Humvee 1151 vehicles"
$ws.Cells.Item(3, 2).Value = 2
$ws.Cells.Item(3, 3).Value = "-"

$ws.Cells.Item(4, 1).Value = "Number of airborne personnel"
$ws.Cells.Item(4, 2).Value = 9
$ws.Cells.Item(4, 3).Value = "-"

$ws.Cells.Item(5, 1).Value = "Range"
$ws.Cells.Item(5, 2).Value = 4000000
$ws.Cells.Item(5, 3).Value = "m"

$ws.Cells.Item(6, 1).Value = "Take-off distance"
$ws.Cells.Item(6, 2).Value = 1093
$ws.Cells.Item(6, 3).Value = "m"

$ws.Cells.Item(7, 1).Value = "Landing distance"
$ws.Cells.Item(7, 2).Value = 975
$ws.Cells.Item(7, 3).Value = "m"

$ws.Cells.Item(8, 1).Value = "Cruise altitude"
$ws.Cells.Item(8, 2).Value = 8535
$ws.Cells.Item(8, 3).Value = "m"

$ws.Cells.Item(9, 1).Value = "Cruise velocity"
$ws.Cells.Item(9, 2).Value = 150
$ws.Cells.Item(9, 3).Value = "m/s"

$ws.Cells.Item(10, 1).Value = "Wing aspect ratio"
$ws.Cells.Item(10, 2).Value = 10.1
$ws.Cells.Item(10, 3).Value = "-"

$ws.Cells.Item(11, 1).Value = "Wing root airfoil"
$ws.Cells.Item(11, 2).NumberFormat = "@"
$ws.Cells.Item(11, 2).Value = "64318"
$ws.Cells.Item(11, 3).Value = "-"

$ws.Cells.Item(12, 1).Value = "Wing tip airfoil"
$ws.Cells.Item(12, 2).NumberFormat = "@"
$ws.Cells.Item(12, 2).Value = "64412"
$ws.Cells.Item(12, 3).Value = "-"

$ws.Cells.Item(13, 1).Value = "Number of engines"
$ws.Cells.Item(13, 2).Value = 4
$ws.Cells.Item(13, 3).Value = "-"

$ws.Cells.Item(14, 1).Value = "Wing position"
$ws.Cells.Item(14, 2).Value = 0.4
$ws.Cells.Item(14, 3).Value = "x/length_fuselage"

$ws.Cells.Item(15, 1).Value = "Horizontal tail thickness ratio"
$ws.Cells.Item(15, 2).NumberFormat = "@"
$ws.Cells.Item(15, 2).Value = "0018"
$ws.Cells.Item(15, 3).Value = "t/c_h_root"

$ws.Cells.Item(16, 1).Value = "Vertical tail thickness ratio"
$ws.Cells.Item(16, 2).NumberFormat = "@"
$ws.Cells.Item(16, 2).Value = "0018"
$ws.Cells.Item(16, 3).Value = "t/c_v_root"

$ws.Cells.Item(17, 1).Value = "Cruise angle of attack"
$ws.Cells.Item(17, 2).Value = 2
$ws.Cells.Item(17, 3).Value = "deg"

$ws.Cells.Item(18, 1).Value = "Ultimate design load factor"
$ws.Cells.Item(18, 2).Value = 3
$ws.Cells.Item(18, 3).Value = "-"

$ws.Cells.Item(19, 1).Value = "Number of fuel tanks"
$ws.Cells.Item(19, 2).Value = 2
$ws.Cells.Item(19, 3).Value = "-"

$ws.Cells.Item(20, 1).Value = "Design propulsive efficiency"
$ws.Cells.Item(20, 2).Value = 0.82
$ws.Cells.Item(20, 3).Value = "-"

# Section header rows (22, 27, 35): label in A, B/C intentionally blank
$ws.Cells.Item(22, 1).Value = "Class I weight estimation"
$ws.Cells.Item(27, 1).Value = "Class II weight estimation"
$ws.Cells.Item(35, 1).Value = "Longitudinal Static Stability"

# Data rows for Class I / Class II / stability sections
$ws.Cells.Item(23, 1).Value = "Operative Empty Weight (OEW)"
$ws.Cells.Item(23, 2).Value = 21093.50425273008
$ws.Cells.Item(23, 3).Value = "kg"

$ws.Cells.Item(24, 1).Value = "Take-Off Weight (TOW)"
$ws.Cells.Item(24, 2).Value = 67707.11367456056
$ws.Cells.Item(24, 3).Value = "kg"

$ws.Cells.Item(25, 1).Value = "Fuel weight (Wf)"
$ws.Cells.Item(25, 2).Value = 35189.60942183048
$ws.Cells.Item(25, 3).Value = "kg"

$ws.Cells.Item(28, 1).Value = "Wing weight"
$ws.Cells.Item(28, 2).Value = 1598.656536988379
$ws.Cells.Item(28, 3).Value = "kg"

$ws.Cells.Item(29, 1).Value = "Fuselage weight"
$ws.Cells.Item(29, 2).Value = 10749.72838332894
$ws.Cells.Item(29, 3).Value = "kg"

$ws.Cells.Item(30, 1).Value = "Engine weight"
$ws.Cells.Item(30, 2).Value = 6243.32026452595
$ws.Cells.Item(30, 3).Value = "kg"

$ws.Cells.Item(31, 1).Value = "Horizontal tail weight"
$ws.Cells.Item(31, 2).Value = 1929.9998132356
$ws.Cells.Item(31, 3).Value = "kg"

$ws.Cells.Item(32, 1).Value = "Vertical tail weight"
$ws.Cells.Item(32, 2).Value = 101.7482430393286
$ws.Cells.Item(32, 3).Value = "kg"

$ws.Cells.Item(33, 1).Value = "Fuel tank weight"
$ws.Cells.Item(33, 2).Value = 470.0510116118821
$ws.Cells.Item(33, 3).Value = "kg"

$ws.Cells.Item(36, 1).Value = "Tailless center of gravity"
$ws.Cells.Item(36, 2).Value = 10.1108299875289
$ws.Cells.Item(36, 3).Value = "m"

$ws.Cells.Item(37, 1).Value = "Total center of gravity"
$ws.Cells.Item(37, 2).Value = 11.15979058353563
$ws.Cells.Item(37, 3).Value = "m"

# Row 39: trailing lone label cell
$ws.Cells.Item(39, 1).Value = "Class I weight estimation"

